$d = $word.ActiveDocument

$replacements = @(
    @("80×52=", "32×87="),
    @("86×44=", "39×67="),
    @("22×21=", "61×51="),
    @("91×88=", "73×47="),
    @("48×15=", "17×46="),
    @("26×19=", "75×27="),
    @("13×34=", "43×61="),
    @("93×21=", "64×85="),
    @("53×99=", "20×21="),
    @("33×11=", "45×23="),
    @("50×21=", "83×76="),
    @("19×49=", "56×21="),
    @("20×49=", "55×12="),
    @("17×86=", "62×41="),
    @("82×11=", "49×73="),
    @("51×82=", "85×15="),
    @("17×81=", "50×30="),
    @("80×32=", "16×38="),
    @("42×73=", "82×94="),
    @("86×82=", "45×20="),
    @("62×14=", "49×84="),
    @("91×32=", "24×88="),
    @("24×73=", "52×48="),
    @("42×25=", "23×54="),
    @("86×72=", "96×47=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
